{"js": "// Fix typo in date: \"May 2023\" -> \"May 2024\"\n// (the CV's \"Scientific Developer ... May 2023 - present\" entry had the\n// wrong year; bump it by one to May 2024, keeping all formatting intact.)\n\nconst body = context.document.body;\n\nconst searchResults = body.search(\"May 2023\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length !== 1) {\n  throw new Error(\n    `Expected exactly one occurrence of \"May 2023\", found ${searchResults.items.length}`\n  );\n}\n\nconst match = searchResults.items[0];\n\n// Replace the matched text in place so the run keeps its original\n// formatting (font, size, etc.) - only the visible characters change.\nmatch.insertText(\"May 2024\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Fix typo in date: \"May 2023\" -> \"May 2024\"\n# (the CV's \"Scientific Developer ... May 2023 - present\" entry had the\n# wrong year; bump it by one to May 2024, keeping all formatting intact.)\n\n$d = $word.ActiveDocument\n\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute(\"May 2023\", $false, $false, $false, $false, $false, $true)\n\nif (-not $found) {\n    throw \"Could not find 'May 2023' in the document\"\n}\n\n# Replacing the Range's Text in place preserves the run's existing\n# formatting (font, size, etc.) - only the visible characters change.\n$searchRange.Text = \"May 2024\"\n"}
